$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

$ws.Range("E2").Value = "2026-02-26 04:48:13"
$ws.Range("E3").Value = "2026-02-26 04:48:15"
$ws.Range("K3").Value = "-0.1 MJ/m2"
$ws.Range("N3").Value = "0.4 °C 4:13 TU"
$ws.Range("O3").Value = "1.3 °C"
$ws.Range("E4").Value = "2026-02-26 04:48:17"
$ws.Range("J4").Value = "1026.4 hPa"
$ws.Range("O4").Value = "7.0 °C"
$ws.Range("E5").Value = "2026-02-26 04:48:20"
$ws.Range("N5").Value = "2.8 °C 4:19 TU"
$ws.Range("E6").Value = "2026-02-26 04:48:22"
$ws.Range("I6").Value = "0.1 mm"
$ws.Range("J6").Value = "1026.1 hPa"
$ws.Range("O6").Value = "9.4 °C"
$ws.Range("E7").Value = "2026-02-26 04:48:25"
$ws.Range("N7").Value = "11.1 °C 4:00 TU"
$ws.Range("O7").Value = "11.7 °C"
$ws.Range("E8").Value = "2026-02-26 04:48:27"
$ws.Range("H8").Value = "'94%"
$ws.Range("M8").Value = "9.7 °C 4:02 TU"
$ws.Range("E9").Value = "2026-02-26 04:48:29"
$ws.Range("N9").Value = "9.3 °C 4:29 TU"
$ws.Range("O9").Value = "10.7 °C"
$ws.Range("E10").Value = "2026-02-26 04:48:32"
$ws.Range("L10").Value = "5.0 km/h - 92º 4:23 TU"
$ws.Range("E11").Value = "2026-02-26 04:48:34"
$ws.Range("H11").Value = "'94%"
$ws.Range("N11").Value = "1.1 °C 4:03 TU"
$ws.Range("O11").Value = "2.1 °C"
$ws.Range("E12").Value = "2026-02-26 04:48:36"
$ws.Range("O12").Value = "9.6 °C"
$ws.Range("E13").Value = "2026-02-26 04:48:38"
$ws.Range("J13").Value = "1031.7 hPa"
$ws.Range("N13").Value = "-2.8 °C 4:29 TU"
$ws.Range("O13").Value = "-0.8 °C"
$ws.Range("E14").Value = "2026-02-26 04:48:41"
$ws.Range("N14").Value = "8.8 °C 4:26 TU"
$ws.Range("O14").Value = "9.8 °C"
$ws.Range("E15").Value = "2026-02-26 04:48:43"
$ws.Range("H15").Value = "'96%"
$ws.Range("N15").Value = "9.0 °C 4:24 TU"
$ws.Range("O15").Value = "10.3 °C"
$ws.Range("E16").Value = "2026-02-26 04:48:45"
$ws.Range("L16").Value = "24.5 km/h - 296º 4:29 TU"
$ws.Range("O16").Value = "1.8 °C"
$ws.Range("E17").Value = "2026-02-26 04:48:47"
$ws.Range("E18").Value = "2026-02-26 04:48:49"
$ws.Range("J18").Value = "1026.4 hPa"
$ws.Range("L18").Value = "5.4 km/h - 264º 4:16 TU"
$ws.Range("N18").Value = "7.4 °C 4:29 TU"
$ws.Range("O18").Value = "8.4 °C"
$ws.Range("E19").Value = "2026-02-26 04:48:51"
$ws.Range("H19").Value = "'67%"
$ws.Range("M19").Value = "8.1 °C 4:01 TU"
$ws.Range("E20").Value = "2026-02-26 04:48:54"
$ws.Range("H20").Value = "'59%"
$ws.Range("E21").Value = "2026-02-26 04:48:56"
$ws.Range("H21").Value = "'84%"
$ws.Range("J21").Value = "1028.7 hPa"
$ws.Range("N21").Value = "2.4 °C 4:23 TU"
$ws.Range("O21").Value = "4.3 °C"
$ws.Range("E22").Value = "2026-02-26 04:48:58"
$ws.Range("H22").Value = "'58%"
$ws.Range("M22").Value = "1.0 °C 4:29 TU"
$ws.Range("E23").Value = "2026-02-26 04:49:01"
$ws.Range("K23").Value = "-0.1 MJ/m2"
$ws.Range("E24").Value = "2026-02-26 04:49:03"
$ws.Range("H24").Value = "'81%"
$ws.Range("J24").Value = "1025.9 hPa"
$ws.Range("N24").Value = "2.6 °C 4:29 TU"
$ws.Range("O24").Value = "7.0 °C"
$ws.Range("E25").Value = "2026-02-26 04:49:05"
$ws.Range("E26").Value = "2026-02-26 04:49:08"
$ws.Range("G26").Value = "1 cm"
$ws.Range("H26").Value = "'44%"
$ws.Range("J26").Value = "1024.9 hPa"
$ws.Range("L26").Value = "15.8 km/h - 351º 4:11 TU"
$ws.Range("E27").Value = "2026-02-26 04:49:10"
$ws.Range("H27").Value = "'57%"
$ws.Range("O27").Value = "2.3 °C"
$ws.Range("E28").Value = "2026-02-26 04:49:13"
$ws.Range("H28").Value = "'88%"
$ws.Range("J28").Value = "1026.1 hPa"
$ws.Range("N28").Value = "7.1 °C 4:28 TU"
$ws.Range("O28").Value = "8.4 °C"
$ws.Range("E29").Value = "2026-02-26 04:49:15"
$ws.Range("E30").Value = "2026-02-26 04:49:17"
$ws.Range("J30").Value = "1026.0 hPa"
$ws.Range("N30").Value = "10.1 °C 4:18 TU"
$ws.Range("O30").Value = "10.8 °C"
$ws.Range("E31").Value = "2026-02-26 04:49:20"
$ws.Range("N31").Value = "9.9 °C 4:05 TU"
$ws.Range("O31").Value = "10.5 °C"
$ws.Range("E32").Value = "2026-02-26 04:49:22"
$ws.Range("H32").Value = "'75%"
$ws.Range("N32").Value = "-0.5 °C 4:23 TU"
$ws.Range("O32").Value = "1.3 °C"
$ws.Range("E33").Value = "2026-02-26 04:49:24"
$ws.Range("J33").Value = "1029.3 hPa"
$ws.Range("N33").Value = "1.3 °C 4:20 TU"
$ws.Range("O33").Value = "2.6 °C"
$ws.Range("E34").Value = "2026-02-26 04:49:27"
$ws.Range("H34").Value = "'50%"
$ws.Range("L34").Value = "25.9 km/h - 39º 4:08 TU"
$ws.Range("O34").Value = "2.7 °C"
$ws.Range("E35").Value = "2026-02-26 04:49:29"
$ws.Range("J35").Value = "1025.3 hPa"
$ws.Range("O35").Value = "8.9 °C"
$ws.Range("E36").Value = "2026-02-26 04:49:31"
$ws.Range("E37").Value = "2026-02-26 04:49:34"
$ws.Range("N37").Value = "1.7 °C 4:29 TU"
$ws.Range("E38").Value = "2026-02-26 04:49:36"
$ws.Range("N38").Value = "5.9 °C 4:04 TU"
$ws.Range("O38").Value = "7.8 °C"
$ws.Range("E39").Value = "2026-02-26 04:49:38"
$ws.Range("O39").Value = "2.8 °C"
$ws.Range("E40").Value = "2026-02-26 04:49:41"
$ws.Range("H40").Value = "'98%"
$ws.Range("J40").Value = "1029.5 hPa"
$ws.Range("N40").Value = "1.3 °C 4:06 TU"
$ws.Range("O40").Value = "2.3 °C"
$ws.Range("E41").Value = "2026-02-26 04:49:43"
$ws.Range("H41").Value = "'98%"
$ws.Range("N41").Value = "6.6 °C 4:26 TU"
$ws.Range("O41").Value = "8.0 °C"
$ws.Range("E42").Value = "2026-02-26 04:49:45"
$ws.Range("O42").Value = "8.6 °C"
$ws.Range("E43").Value = "2026-02-26 04:49:48"
$ws.Range("L43").Value = "5.8 km/h - 213º 4:16 TU"
$ws.Range("O43").Value = "2.9 °C"
$ws.Range("E44").Value = "2026-02-26 04:49:50"
$ws.Range("H44").Value = "'60%"
$ws.Range("L44").Value = "18.0 km/h - 69º 4:25 TU"
$ws.Range("M44").Value = "1.6 °C 4:26 TU"
$ws.Range("O44").Value = "0.0 °C"
$ws.Range("E45").Value = "2026-02-26 04:49:52"
$ws.Range("J45").Value = "1027.3 hPa"
$ws.Range("O45").Value = "6.2 °C"
$ws.Range("E46").Value = "2026-02-26 04:49:54"
$ws.Range("J46").Value = "1025.9 hPa"
$ws.Range("N46").Value = "5.5 °C 4:29 TU"
$ws.Range("O46").Value = "7.6 °C"
